# EFO update - fix wrong lymph/myelo entries
# Row 3 (Disease Ontology) keeps version "v2023-12-20"
# Row 4 (Experimental Factor Ontology) version bumped to "v3.62.0"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phen_oncox")

$ws.Range("E3").Value = "v2023-12-20"
$ws.Range("E4").Value = "v3.62.0"

# Move the active selection from E3 to E4
$ws.Range("E4").Select()
